$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 46509.637
$ws.Range("I28").Value2 = 63846.875
$ws.Range("J28").Value2 = 277
$ws.Range("K28").Value2 = 63846.875
$ws.Range("L28").Value2 = 277
$ws.Range("M28").Value2 = -63361.875
$ws.Range("N28").Value2 = -1247
$ws.Range("H96").Value2 = 758.36365
$ws.Range("I96").Value2 = 680.25
$ws.Range("K96").Value2 = 2040.75
$ws.Range("M96").Value2 = -667.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1617.909
$ws.Range("I2").Value2 = 1655.2222
$ws.Range("K2").Value2 = 1655.2222
$ws.Range("M2").Value2 = -1542.2222
$ws.Range("H32").Value2 = 4135.024
$ws.Range("I32").Value2 = 4135.024
$ws.Range("K32").Value2 = 4135.024
$ws.Range("M32").Value2 = -3848.024
$ws.Range("H45").Value2 = 2530.6667
$ws.Range("I45").Value2 = 1539.4286
$ws.Range("K45").Value2 = 1539.4286
$ws.Range("M45").Value2 = -1162.4286
$ws.Range("H74").Value2 = 1221.0869
$ws.Range("I74").Value2 = 1221.0869
$ws.Range("K74").Value2 = 1221.0869
$ws.Range("M74").Value2 = -347.0869
$ws.Range("H77").Value2 = 1221.0869
$ws.Range("I77").Value2 = 1221.0869
$ws.Range("K77").Value2 = 6105.4345
$ws.Range("M77").Value2 = -1737.4345
$ws.Range("H102").Value2 = 1694
$ws.Range("I102").Value2 = 1651.5
$ws.Range("K102").Value2 = 1651.5
$ws.Range("M102").Value2 = -29.5
$ws.Range("H110").Value2 = 296246.34
$ws.Range("I110").Value2 = 359512.5
$ws.Range("J110").Value2 = 1004.3333
$ws.Range("K110").Value2 = 359512.5
$ws.Range("L110").Value2 = 1004.3333
$ws.Range("M110").Value2 = -357467.5
$ws.Range("N110").Value2 = -5094.3333
$ws.Range("H116").Value2 = 1617.909
$ws.Range("I116").Value2 = 1655.2222
$ws.Range("K116").Value2 = 1655.2222
$ws.Range("M116").Value2 = 638.7778000000001
$ws.Range("H132").Value2 = 2716.4048
$ws.Range("I132").Value2 = 2897.7878
$ws.Range("J132").Value2 = 2051.3333
$ws.Range("K132").Value2 = 8693.3634
$ws.Range("L132").Value2 = 6153.999899999999
$ws.Range("M132").Value2 = -6163.3634
$ws.Range("N132").Value2 = -11213.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1617.909
$ws.Range("I3").Value2 = 1655.2222
$ws.Range("K3").Value2 = 1655.2222
$ws.Range("M3").Value2 = -1541.2222
$ws.Range("H82").Value2 = 45166.668
$ws.Range("I82").Value2 = 7750
$ws.Range("K82").Value2 = 7750
$ws.Range("M82").Value2 = -7367
$ws.Range("H85").Value2 = 45166.668
$ws.Range("I85").Value2 = 7750
$ws.Range("K85").Value2 = 7750
$ws.Range("M85").Value2 = -6424
$ws.Range("H86").Value2 = 1373
$ws.Range("I86").Value2 = 1398.5
$ws.Range("K86").Value2 = 1398.5
$ws.Range("M86").Value2 = -275.5
$ws.Range("H89").Value2 = 1373
$ws.Range("I89").Value2 = 1398.5
$ws.Range("K89").Value2 = 6992.5
$ws.Range("M89").Value2 = -1376.5
$ws.Range("H94").Value2 = 0
$ws.Range("I94").Value2 = 0
$ws.Range("J94").Value2 = 0
$ws.Range("K94").Value2 = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value2 = 0
$ws.Range("H107").Value2 = 1916.4615
$ws.Range("I107").Value2 = 1460.8572
$ws.Range("J107").Value2 = 2448
$ws.Range("K107").Value2 = 1460.8572
$ws.Range("L107").Value2 = 2448
$ws.Range("M107").Value2 = 459.1428000000001
$ws.Range("N107").Value2 = -6288

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 6690.8823
$ws.Range("I16").Value2 = 2372.182
$ws.Range("J16").Value2 = 14608.5
$ws.Range("K16").Value2 = 2372.182
$ws.Range("L16").Value2 = 14608.5
$ws.Range("M16").Value2 = -2085.182
$ws.Range("N16").Value2 = -15182.5
$ws.Range("H58").Value2 = 4853.911
$ws.Range("I58").Value2 = 4319.839
$ws.Range("J58").Value2 = 6036.5
$ws.Range("K58").Value2 = 4319.839
$ws.Range("L58").Value2 = 6036.5
$ws.Range("M58").Value2 = -4116.839
$ws.Range("N58").Value2 = -6442.5
$ws.Range("H113").Value2 = 6690.8823
$ws.Range("I113").Value2 = 2372.182
$ws.Range("J113").Value2 = 14608.5
$ws.Range("K113").Value2 = 2372.182
$ws.Range("L113").Value2 = 14608.5
$ws.Range("M113").Value2 = -202.1819999999998
$ws.Range("N113").Value2 = -18948.5
$ws.Range("H132").Value2 = 1252.7778
$ws.Range("I132").Value2 = 1096.875
$ws.Range("J132").Value2 = 2500
$ws.Range("K132").Value2 = 3290.625
$ws.Range("L132").Value2 = 7500
$ws.Range("M132").Value2 = -760.625
$ws.Range("N132").Value2 = -12560
$ws.Range("H134").Value2 = 219248.66
$ws.Range("I134").Value2 = 1925.8636
$ws.Range("K134").Value2 = 5777.5908
$ws.Range("M134").Value2 = -3242.5908
$ws.Range("H136").Value2 = 4853.911
$ws.Range("I136").Value2 = 4319.839
$ws.Range("J136").Value2 = 6036.5
$ws.Range("K136").Value2 = 12959.517
$ws.Range("L136").Value2 = 18109.5
$ws.Range("M136").Value2 = -10409.517
$ws.Range("N136").Value2 = -23209.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value2 = 166879.22
$ws.Range("J81").Value2 = 166879.22
$ws.Range("L81").Value2 = 500637.66
$ws.Range("N81").Value2 = -502883.66
$ws.Range("H84").Value2 = 166879.22
$ws.Range("J84").Value2 = 166879.22
$ws.Range("L84").Value2 = 1501912.98
$ws.Range("N84").Value2 = -1513144.98
$ws.Range("H137").Value2 = 3015
$ws.Range("I137").Value2 = 1765.6666
$ws.Range("J137").Value2 = 7700
$ws.Range("K137").Value2 = 5296.9998
$ws.Range("L137").Value2 = 23100
$ws.Range("M137").Value2 = -196.9997999999996
$ws.Range("N137").Value2 = -33300

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 1729.7916
$ws.Range("I102").Value2 = 639.7778
$ws.Range("J102").Value2 = 4999.8335
$ws.Range("K102").Value2 = 639.7778
$ws.Range("L102").Value2 = 4999.8335
$ws.Range("M102").Value2 = 982.2222
$ws.Range("N102").Value2 = -8243.833500000001
$ws.Range("H113").Value2 = 462125.6
$ws.Range("I113").Value2 = 1430485.8
$ws.Range("J113").Value2 = 10224.2
$ws.Range("K113").Value2 = 1430485.8
$ws.Range("L113").Value2 = 10224.2
$ws.Range("M113").Value2 = -1428315.8
$ws.Range("N113").Value2 = -14564.2
$ws.Range("H126").Value2 = 4106.091
$ws.Range("J126").Value2 = 4174.25
$ws.Range("L126").Value2 = 12522.75
$ws.Range("N126").Value2 = -17462.75
$ws.Range("H132").Value2 = 42403.348
$ws.Range("I132").Value2 = 4604.579
$ws.Range("K132").Value2 = 13813.737
$ws.Range("M132").Value2 = -11283.737

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 250001250
$ws.Range("I16").Value2 = 500000740
$ws.Range("K16").Value2 = 500000740
$ws.Range("M16").Value2 = -500000570
$ws.Range("H40").Value2 = 149068
$ws.Range("I40").Value2 = 254869
$ws.Range("K40").Value2 = 254869
$ws.Range("M40").Value2 = -254733
$ws.Range("H61").Value2 = 5654.5454
$ws.Range("I61").Value2 = 5900.75
$ws.Range("J61").Value2 = 4998
$ws.Range("K61").Value2 = 5900.75
$ws.Range("L61").Value2 = 4998
$ws.Range("M61").Value2 = -5698.75
$ws.Range("N61").Value2 = -5402
$ws.Range("H93").Value2 = 4615.3335
$ws.Range("I93").Value2 = 4498.75
$ws.Range("K93").Value2 = 4498.75
$ws.Range("M93").Value2 = -3250.75
$ws.Range("H113").Value2 = 5654.5454
$ws.Range("I113").Value2 = 5900.75
$ws.Range("J113").Value2 = 4998
$ws.Range("K113").Value2 = 5900.75
$ws.Range("L113").Value2 = 4998
$ws.Range("M113").Value2 = -3730.75
$ws.Range("N113").Value2 = -9338
$ws.Range("H122").Value2 = 3334232
$ws.Range("J122").Value2 = 5000558
$ws.Range("L122").Value2 = 15001674
$ws.Range("N122").Value2 = -15006574
$ws.Range("H132").Value2 = 2884.4285
$ws.Range("I132").Value2 = 2884.4285
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 8653.2855
$ws.Range("L132").Value2 = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value2 = -6123.2855
$ws.Range("H136").Value2 = 915522.1
$ws.Range("I136").Value2 = 1339111.6
$ws.Range("J136").Value2 = 7830.2856
$ws.Range("K136").Value2 = 4017334.8
$ws.Range("L136").Value2 = 23490.8568
$ws.Range("M136").Value2 = -4014784.8
$ws.Range("N136").Value2 = -28590.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value2 = 1500
$ws.Range("I14").Value2 = 1750
$ws.Range("J14").Value2 = 1000
$ws.Range("K14").Value2 = 1750
$ws.Range("L14").Value2 = 1000
$ws.Range("M14").Value2 = -1582
$ws.Range("N14").Value2 = -1336
$ws.Range("H17").Value2 = 5000
$ws.Range("I17").Value2 = 5000
$ws.Range("K17").Value2 = 5000
$ws.Range("M17").Value2 = -4828
$ws.Range("H100").Value2 = 740.9
$ws.Range("I100").Value2 = 740.9
$ws.Range("K100").Value2 = 1481.8
$ws.Range("M100").Value2 = -940.8
$ws.Range("H107").Value2 = 719.9286
$ws.Range("I107").Value2 = 781.0952
$ws.Range("J107").Value2 = 536.4286
$ws.Range("K107").Value2 = 2343.2856
$ws.Range("L107").Value2 = 1609.2858
$ws.Range("M107").Value2 = -423.2856000000002
$ws.Range("N107").Value2 = -5449.2858
$ws.Range("H113").Value2 = 812.25
$ws.Range("I113").Value2 = 749.6667
$ws.Range("J113").Value2 = 1000
$ws.Range("K113").Value2 = 2249.0001
$ws.Range("L113").Value2 = 3000
$ws.Range("M113").Value2 = -79.0001000000002
$ws.Range("N113").Value2 = -7340
$ws.Range("H126").Value2 = 1670.625
$ws.Range("I126").Value2 = 1670.625
$ws.Range("K126").Value2 = 5011.875
$ws.Range("M126").Value2 = -2541.875
$ws.Range("H136").Value2 = 8424761
$ws.Range("I136").Value2 = 11077522
$ws.Range("J136").Value2 = 201199
$ws.Range("K136").Value2 = 33232566
$ws.Range("L136").Value2 = 603597
$ws.Range("M136").Value2 = -33230016
$ws.Range("N136").Value2 = -608697
